$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column D, shifting the existing
# D:K data block to E:L. Excel's default Insert() copies formatting
# from the column to the left, so we restore per-cell formatting by
# copying it back from the (now shifted) column E, which carries the
# original column D formatting. Only the three data blocks (the rows
# that actually carry figures in column D) are touched, so the blank
# spacer rows (36, 78) and the section-header rows (37, 79, which only
# have a value in column B) are left untouched.
$ws.Columns("D:D").Insert()
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the newly inserted column D with the new period's figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 7270400
$ws.Range("D9").Value = 6130000
$ws.Range("D10").Value = 1140400
$ws.Range("D12").Value = 146200
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 552700
$ws.Range("D15").Value = 99400
$ws.Range("D17").Value = 7167800
$ws.Range("D18").Value = 102600
$ws.Range("D20").Value = -200
$ws.Range("D21").Value = 631200
$ws.Range("D22").Value = 216300
$ws.Range("D23").Value = -113900
$ws.Range("D24").Value = -57100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -56800
$ws.Range("D27").Value = -57500
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 200
$ws.Range("D33").Value = -57500
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -57500
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 476400
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 966500
$ws.Range("D44").Value = 459700
$ws.Range("D45").Value = 127200
$ws.Range("D46").Value = 2029800
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 2514400
$ws.Range("D49").Value = 2252900
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 713600
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 7510700
$ws.Range("D57").Value = 840200
$ws.Range("D58").Value = 121600
$ws.Range("D59").Value = 395000
$ws.Range("D60").Value = 1356800
$ws.Range("D61").Value = 3686800
$ws.Range("D62").Value = 980800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 6026800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 703500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1483900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -57500
$ws.Range("D83").Value = 528800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 771500
$ws.Range("D91").Value = -525200
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -478200
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -184500
$ws.Range("D101").Value = -6700
$ws.Range("D102").Value = 102100

